$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 575.3333
$ws.Range("J29").Value = 712.5
$ws.Range("L29").Value = 2137.5
$ws.Range("N29").Value = -2699.5
$ws.Range("H38").Value = 2481482.8
$ws.Range("I38").Value = 2481482.8
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 7444448.399999999
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -7444076.399999999
$ws.Range("N38").ClearContents()
$ws.Range("H58").Value = 1516204
$ws.Range("I58").Value = 2164720
$ws.Range("J58").Value = 3000
$ws.Range("K58").Value = 6494160
$ws.Range("L58").Value = 9000
$ws.Range("M58").Value = -6494010
$ws.Range("N58").Value = -9300
$ws.Range("H74").Value = 4431.35
$ws.Range("I74").Value = 4384.75
$ws.Range("J74").Value = 4462.4165
$ws.Range("K74").Value = 4384.75
$ws.Range("L74").Value = 4462.4165
$ws.Range("M74").Value = -3448.75
$ws.Range("N74").Value = -6334.4165
$ws.Range("H77").Value = 4431.35
$ws.Range("I77").Value = 4384.75
$ws.Range("J77").Value = 4462.4165
$ws.Range("K77").Value = 21923.75
$ws.Range("L77").Value = 22312.0825
$ws.Range("M77").Value = -17243.75
$ws.Range("N77").Value = -31672.0825
$ws.Range("H80").Value = 67347.664
$ws.Range("I80").Value = 374
$ws.Range("J80").Value = 100834.5
$ws.Range("K80").Value = 1122
$ws.Range("L80").Value = 302503.5
$ws.Range("M80").Value = -124
$ws.Range("N80").Value = -304499.5
$ws.Range("H83").Value = 67347.664
$ws.Range("I83").Value = 374
$ws.Range("J83").Value = 100834.5
$ws.Range("K83").Value = 3366
$ws.Range("L83").Value = 907510.5
$ws.Range("M83").Value = 1626
$ws.Range("N83").Value = -917494.5
$ws.Range("H86").Value = 6987.643
$ws.Range("I86").Value = 5943.75
$ws.Range("J86").Value = 8379.5
$ws.Range("K86").Value = 5943.75
$ws.Range("L86").Value = 8379.5
$ws.Range("M86").Value = -4820.75
$ws.Range("N86").Value = -10625.5
$ws.Range("H89").Value = 6987.643
$ws.Range("I89").Value = 5943.75
$ws.Range("J89").Value = 8379.5
$ws.Range("K89").Value = 29718.75
$ws.Range("L89").Value = 41897.5
$ws.Range("M89").Value = -24102.75
$ws.Range("N89").Value = -53129.5
$ws.Range("H135").Value = 976.1
$ws.Range("I135").Value = 1065.3572
$ws.Range("J135").Value = 767.8333
$ws.Range("K135").Value = 9588.2148
$ws.Range("L135").Value = 6910.4997
$ws.Range("M135").Value = -7053.2148
$ws.Range("N135").Value = -11980.4997
$ws.Range("H138").Value = 2819.1475
$ws.Range("I138").Value = 1809.7646
$ws.Range("J138").Value = 3209.1365
$ws.Range("K138").Value = 5429.293799999999
$ws.Range("L138").Value = 9627.4095
$ws.Range("M138").Value = -289.2937999999995
$ws.Range("N138").Value = -19907.4095

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H5").Value = 151.9375
$ws.Range("I5").Value = 123.1
$ws.Range("K5").Value = 123.1
$ws.Range("M5").Value = -11.09999999999999

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 151.9375
$ws.Range("I4").Value = 123.1
$ws.Range("K4").Value = 123.1
$ws.Range("M4").Value = -8.099999999999994
$ws.Range("H134").Value = 1485.0488
$ws.Range("I134").Value = 1351.3243
$ws.Range("K134").Value = 4053.9729
$ws.Range("M134").Value = -1518.9729

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 145.41667
$ws.Range("I7").Value = 50.692307
$ws.Range("K7").Value = 50.692307
$ws.Range("M7").Value = 62.307693
$ws.Range("H31").Value = 1366.5278
$ws.Range("I31").Value = 850.25
$ws.Range("J31").Value = 2708.85
$ws.Range("K31").Value = 850.25
$ws.Range("L31").Value = 2708.85
$ws.Range("M31").Value = -555.25
$ws.Range("N31").Value = -3298.85
$ws.Range("H34").Value = 1366.5278
$ws.Range("I34").Value = 850.25
$ws.Range("J34").Value = 2708.85
$ws.Range("K34").Value = 850.25
$ws.Range("L34").Value = 2708.85
$ws.Range("M34").Value = -648.25
$ws.Range("N34").Value = -3112.85
$ws.Range("H58").Value = 10742.223
$ws.Range("I58").Value = 1228.8
$ws.Range("J58").Value = 37923.43
$ws.Range("K58").Value = 1228.8
$ws.Range("L58").Value = 37923.43
$ws.Range("M58").Value = -1025.8
$ws.Range("N58").Value = -38329.43
$ws.Range("H99").Value = 8072.375
$ws.Range("I99").Value = 1781.375
$ws.Range("J99").Value = 14363.375
$ws.Range("K99").Value = 1781.375
$ws.Range("L99").Value = 14363.375
$ws.Range("M99").Value = -283.375
$ws.Range("N99").Value = -17359.375
$ws.Range("H107").Value = 922.9231
$ws.Range("I107").Value = 1274.8334
$ws.Range("J107").Value = 621.2857
$ws.Range("K107").Value = 1274.8334
$ws.Range("L107").Value = 621.2857
$ws.Range("M107").Value = 645.1666
$ws.Range("N107").Value = -4461.2857
$ws.Range("H126").Value = 8072.375
$ws.Range("I126").Value = 1781.375
$ws.Range("J126").Value = 14363.375
$ws.Range("K126").Value = 5344.125
$ws.Range("L126").Value = 43090.125
$ws.Range("M126").Value = -2874.125
$ws.Range("N126").Value = -48030.125
$ws.Range("H136").Value = 10742.223
$ws.Range("I136").Value = 1228.8
$ws.Range("J136").Value = 37923.43
$ws.Range("K136").Value = 3686.4
$ws.Range("L136").Value = 113770.29
$ws.Range("M136").Value = -1136.4
$ws.Range("N136").Value = -118870.29

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 1778.2
$ws.Range("I23").Value = 3433.3333
$ws.Range("J23").Value = 1068.8572
$ws.Range("K23").Value = 10299.9999
$ws.Range("L23").Value = 3206.5716
$ws.Range("M23").Value = -10064.9999
$ws.Range("N23").Value = -3676.5716
$ws.Range("H68").Value = 17060.842
$ws.Range("I68").Value = 929.25
$ws.Range("J68").Value = 24563.906
$ws.Range("K68").Value = 2787.75
$ws.Range("L68").Value = 73691.71799999999
$ws.Range("M68").Value = -1976.75
$ws.Range("N68").Value = -75313.71799999999
$ws.Range("H71").Value = 17060.842
$ws.Range("I71").Value = 929.25
$ws.Range("J71").Value = 24563.906
$ws.Range("K71").Value = 8363.25
$ws.Range("L71").Value = 221075.154
$ws.Range("M71").Value = -4307.25
$ws.Range("N71").Value = -229187.154
$ws.Range("H75").Value = 3999
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 3999
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 11997
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -13993
$ws.Range("H78").Value = 3999
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 3999
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 35991
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -45975
$ws.Range("H81").Value = 16668291
$ws.Range("I81").Value = 910
$ws.Range("J81").Value = 41669364
$ws.Range("K81").Value = 2730
$ws.Range("L81").Value = 125008092
$ws.Range("M81").Value = -1607
$ws.Range("N81").Value = -125010338
$ws.Range("H84").Value = 16668291
$ws.Range("I84").Value = 910
$ws.Range("J84").Value = 41669364
$ws.Range("K84").Value = 8190
$ws.Range("L84").Value = 375024276
$ws.Range("M84").Value = -2574
$ws.Range("N84").Value = -375035508
$ws.Range("H86").Value = 1233.1538
$ws.Range("I86").Value = 1000.3333
$ws.Range("J86").Value = 1432.7142
$ws.Range("K86").Value = 3000.9999
$ws.Range("L86").Value = 4298.142599999999
$ws.Range("M86").Value = -1814.9999
$ws.Range("N86").Value = -6670.142599999999
$ws.Range("H87").Value = 11579.786
$ws.Range("I87").Value = 3188.1428
$ws.Range("K87").Value = 9564.428400000001
$ws.Range("M87").Value = -8316.428400000001
$ws.Range("H89").Value = 1233.1538
$ws.Range("I89").Value = 1000.3333
$ws.Range("J89").Value = 1432.7142
$ws.Range("K89").Value = 9002.9997
$ws.Range("L89").Value = 12894.4278
$ws.Range("M89").Value = -3074.9997
$ws.Range("N89").Value = -24750.4278
$ws.Range("H90").Value = 11579.786
$ws.Range("I90").Value = 3188.1428
$ws.Range("K90").Value = 28693.2852
$ws.Range("M90").Value = -22453.2852
$ws.Range("H131").Value = 758457
$ws.Range("I131").Value = 532.3333
$ws.Range("J131").Value = 785207.25
$ws.Range("K131").Value = 1596.9999
$ws.Range("L131").Value = 2355621.75
$ws.Range("M131").Value = 3443.0001
$ws.Range("N131").Value = -2365701.75
$ws.Range("H137").Value = 30954892
$ws.Range("I137").Value = 2346.6667
$ws.Range("J137").Value = 54169304
$ws.Range("K137").Value = 7040.000100000001
$ws.Range("L137").Value = 162507912
$ws.Range("M137").Value = -1940.000100000001
$ws.Range("N137").Value = -162518112

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 122.94118
$ws.Range("I2").Value = 22.166666
$ws.Range("J2").Value = 177.90909
$ws.Range("K2").Value = 22.166666
$ws.Range("L2").Value = 177.90909
$ws.Range("M2").Value = 90.83333400000001
$ws.Range("N2").Value = -403.90909
$ws.Range("H64").Value = 25666.666
$ws.Range("J64").Value = 25666.666
$ws.Range("L64").Value = 25666.666
$ws.Range("N64").Value = -26162.666
$ws.Range("H67").Value = 25666.666
$ws.Range("J67").Value = 25666.666
$ws.Range("L67").Value = 25666.666
$ws.Range("N67").Value = -27382.666
$ws.Range("H70").Value = 84079.24000000001
$ws.Range("I70").Value = 121634.53
$ws.Range("J70").Value = 4274.25
$ws.Range("K70").Value = 121634.53
$ws.Range("L70").Value = 4274.25
$ws.Range("M70").Value = -121364.53
$ws.Range("N70").Value = -4814.25
$ws.Range("H73").Value = 84079.24000000001
$ws.Range("I73").Value = 121634.53
$ws.Range("J73").Value = 4274.25
$ws.Range("K73").Value = 121634.53
$ws.Range("L73").Value = 4274.25
$ws.Range("M73").Value = -120698.53
$ws.Range("N73").Value = -6146.25

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H119").Value = 38712
$ws.Range("J119").Value = 38712
$ws.Range("L119").Value = 38712
$ws.Range("N119").Value = -48388

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 671.5
$ws.Range("I113").Value = 417.57144
$ws.Range("J113").Value = 869
$ws.Range("K113").Value = 1252.71432
$ws.Range("L113").Value = 2607
$ws.Range("M113").Value = 917.28568
$ws.Range("N113").Value = -6947
$ws.Range("H136").Value = 1305.4762
$ws.Range("I136").Value = 623.7273
$ws.Range("J136").Value = 2055.4
$ws.Range("K136").Value = 1871.1819
$ws.Range("L136").Value = 6166.200000000001
$ws.Range("M136").Value = 678.8181
$ws.Range("N136").Value = -11266.2
